$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "questions" - remove all sample question rows, keep header only
# ---------------------------------------------------------------------------
$wsQuestions = $wb.Worksheets.Item("questions")
$wsQuestions.Rows("2:5").Delete()
$wsQuestions.Range("A1").Select()

# ---------------------------------------------------------------------------
# Sheet 2: "archieves questions" - only the active selection changes
# ---------------------------------------------------------------------------
$wsArchive = $wb.Worksheets.Item("archieves questions")
$wsArchive.Range("A1").Select()

# ---------------------------------------------------------------------------
# Sheet 3: "students" - add leaderboard score column + new students
# ---------------------------------------------------------------------------
$wsStudents = $wb.Worksheets.Item("students")

# New "score" header column
$wsStudents.Range("D1").Value = "score"

# Update the existing student row (id/name changed, score added)
$wsStudents.Range("A2").Value = "vfm4w"
$wsStudents.Range("B2").Value = "Aloa, Mikaella"
$wsStudents.Range("D2").Value = 5

# Additional student rows
$wsStudents.Range("A3").Value = "2n55d"
$wsStudents.Range("B3").Value = "Kimmy, Rheign"
$wsStudents.Range("C3").Value = "8349c112e6b9b83a9296d60d1f7783551ebe7941c2ddab7597ccd727f338bf81"
$wsStudents.Range("D3").Value = 4

$wsStudents.Range("A4").Value = "g1kwu"
$wsStudents.Range("B4").Value = "Rye, Rhianne"
$wsStudents.Range("C4").Value = "8349c112e6b9b83a9296d60d1f7783551ebe7941c2ddab7597ccd727f338bf81"
$wsStudents.Range("D4").Value = 3

$wsStudents.Range("A5").Value = "my6xx"
$wsStudents.Range("B5").Value = "Sesgundo, Ruina"
$wsStudents.Range("C5").Value = "8349c112e6b9b83a9296d60d1f7783551ebe7941c2ddab7597ccd727f338bf81"
$wsStudents.Range("D5").Value = 6

# Column widths for the new leaderboard layout (best-fit style columns)
$wsStudents.Columns.Item(2).ColumnWidth = 15.7109375
$wsStudents.Columns.Item(3).ColumnWidth = 67.42578125

$wsStudents.Range("D6").Select()

# ---------------------------------------------------------------------------
# Sheet 4: "teachers" - rename the existing teacher entry
# ---------------------------------------------------------------------------
$wsTeachers = $wb.Worksheets.Item("teachers")
$wsTeachers.Range("B2").Value = "Malabanan, RySes"
$wsTeachers.Range("B2").Select()

# Restore "teachers" as the active tab/sheet, as it was before editing
$wsTeachers.Activate()
